$d = $word.ActiveDocument

$replacements = @(
    @{old = "196 (71.5)"; new = "217 (72.1)"},
    @{old = "69 (74.2)";  new = "76 (74.5)"},
    @{old = "61 (66.3)";  new = "69 (68.3)"},
    @{old = "66 (74.2)";  new = "72 (73.5)"},
    @{old = "32 (11.7)";  new = "32 (10.6)"},
    @{old = "13 (14.0)";  new = "13 (12.7)"},
    @{old = "11 (12.0)";  new = "11 (10.9)"},
    @{old = "8 (9.0)";    new = "8 (8.2)"},
    @{old = "130 (47.4)"; new = "143 (47.5)"},
    @{old = "43 (46.2)";  new = "47 (46.1)"},
    @{old = "40 (43.5)";  new = "45 (44.6)"},
    @{old = "47 (52.8)";  new = "51 (52.0)"},
    @{old = "62 (22.6)";  new = "68 (22.6)"},
    @{old = "18 (19.4)";  new = "20 (19.6)"},
    @{old = "26 (28.3)";  new = "27 (26.7)"},
    @{old = "18 (20.2)";  new = "21 (21.4)"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
